# Commit 12: limitado o botão de limpar banco de dados aos usuarios administradores
# Data change: a new student ("Aluno 161") was added to the "teste_alunos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new student record as row 82 (directly below the last existing
# data row, 81), following the same column layout as every other row:
#   A = R.A. (registration number)
#   B = Nome do estudante
#   C = Série/turma
#   D = Endereço
#   E = Responsável 1 (pai)
#   F = Responsável 2 (mãe)
#   G = Contato(s)
$ws.Range("A82").Value = 202650
$ws.Range("B82").Value = "Aluno 161"
$ws.Range("C82").Value = "6B"
$ws.Range("D82").Value = "Rua do Aluno 161"
$ws.Range("E82").Value = "Pai do Aluno 161"
$ws.Range("F82").Value = "Mãe do Aluno 161"
$ws.Range("G82").Value = 11987654499

# Match the author's final cursor/selection position in the sheet after the
# edit (they had scrolled down and landed on the empty row right after the
# new entry).
$ws.Range("E87").Select()
